$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update student record in row 2
$ws.Range("A2").Value = "RC3"
$ws.Range("J2").Value = "Mehboobnagar"
$ws.Range("C2").Value = "Ramya"

# Update the selected cell shown in the saved view
$ws.Range("D8").Select()
